$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'52.844.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -12.95%  "

$ws.Range("D3").Value = "'2.313.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -20.35%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'435.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -17.57%  "

$ws.Range("D6").Value = "'120.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -16.52%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "'0.472"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -14.93%  "

$ws.Range("D9").Value = "'2.316.44"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -20.44%  "

$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").Value = "'5.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -13.09%  "

$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.0906"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -16.21%  "

$ws.Range("D12").Value = "'0.307"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -15.13%  "

$ws.Range("E13").Value = "  -3.90%  "

$ws.Range("D14").Value = "'2.727.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -20.04%  "

$ws.Range("D15").Value = "'52.838.06"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -12.90%  "

$ws.Range("D16").Value = "'19.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -16.76%  "

$ws.Range("D17").Value = "'0.0000119"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -16.18%  "

$ws.Range("D18").Value = "'2.332.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -19.99%  "

$ws.Range("D19").Value = "'3.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -21.66%  "

$ws.Range("D20").Value = "'299.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -17.11%  "

$ws.Range("D21").Value = "'9.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -22.58%  "

$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.12%  "

$ws.Range("E23").Value = "  -2.21%  "

$ws.Range("D24").Value = "'5.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -19.56%  "

$ws.Range("D25").Value = "'55.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -14.47%  "

$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.44%  "

$ws.Range("E27").Value = "  -15.13%  "

$ws.Range("D28").Value = "'0.366"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -19.60%  "

$ws.Range("D29").Value = "'6.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -12.12%  "

$ws.Range("E30").Value = "  -0.20%  "

$ws.Range("D31").Value = "'0.0₃0695"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -18.59%  "

$ws.Range("D32").Value = "'142.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.88%  "

$ws.Range("D33").Value = "'17.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -13.49%  "

$ws.Range("D34").Value = "'1.34"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -20.40%  "

$ws.Range("D35").Value = "'4.73"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -15.49%  "

$ws.Range("D36").Value = "'3.52"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -19.35%  "

$ws.Range("D37").Value = "'0.828"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -17.24%  "

$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").Value = "'33.32"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -11.94%  "

$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").Value = "'0.997"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.04%  "

$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "'0.996"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -17.56%  "

$ws.Range("D41").Value = "'10.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.50%  "

$ws.Range("D42").Value = "'3.15"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -15.37%  "

$ws.Range("D43").Value = "'0.0499"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -14.96%  "

$ws.Range("D44").Value = "'1.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -18.99%  "

$ws.Range("D45").Value = "'1.887.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -17.86%  "

$ws.Range("D46").Value = "'0.521"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -19.69%  "

$ws.Range("D47").Value = "'0.0208"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -12.73%  "

$ws.Range("D48").Value = "'0.0831"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -10.11%  "

$ws.Range("D49").Value = "'15.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -23.93%  "

$ws.Range("B50").Value = "ZEEBU"
$ws.Range("C50").Value = "https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu"
$ws.Range("D50").Value = "'4.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.51%  "

$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'3.86"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -22.65%  "
